$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individual")

# --- Fix swapped Fuel/Oxidizer melting & boiling point values, with refined precision ---
$ws.Cells.Item(13, 2).Value = -114
$ws.Cells.Item(13, 3).Value = -218.79
$ws.Cells.Item(14, 2).Value = 78
$ws.Cells.Item(14, 3).Value = -182.96

# --- Units column (E) marking degrees Celsius ---
$ws.Cells.Item(13, 5).Value = "C"
$ws.Cells.Item(14, 5).Value = "C"
$ws.Cells.Item(15, 5).Value = "C"

# --- New "Standard Range" row ---
$ws.Cells.Item(15, 1).Value = "Standard Range"
$ws.Cells.Item(15, 2).Value = '"-114 <> 78"'
$ws.Cells.Item(15, 3).Value = '"-218 <> -182.96"'
$ws.Cells.Item(15, 3).Font.Name = "Arial"
$ws.Cells.Item(15, 3).Font.Size = 12
$ws.Cells.Item(15, 3).Font.Color = 2236962

# Widen column C to fit new content
$ws.Columns.Item(3).ColumnWidth = 17.1640625
